# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Mon Sep  4 20:41:43 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.846.16'
$ws.Range('E2').Value = '  -0.63%  '

$ws.Range('D3').Value = '1.630.33'
$ws.Range('E3').Value = '  -0.70%  '

$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').Value = '''215.56'
$ws.Range('E5').Value = '  +0.38%  '

$ws.Range('E6').Value = '  +0.42%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = '''0.2560'
$ws.Range('E8').Value = '  -0.08%  '

$ws.Range('D9').Value = '''0.06334'
$ws.Range('E9').Value = '  -0.30%  '

$ws.Range('D10').Value = '''19.45'
$ws.Range('E10').Value = '  -0.40%  '

$ws.Range('D11').Value = '''0.07781'
$ws.Range('E11').Value = '  +0.27%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.639.69'
$ws.Range('E12').Value = '  -0.13%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.240'
$ws.Range('E13').Value = '  -0.97%  '

$ws.Range('D14').Value = '1.854.04'
$ws.Range('E14').Value = '  -0.80%  '

$ws.Range('D15').Value = '''0.5519'

$ws.Range('D16').Value = '''63.68'
$ws.Range('E16').Value = '  -0.84%  '

$ws.Range('D17').Value = '0.0₅7605'
$ws.Range('E17').Value = '  -1.38%  '

$ws.Range('D18').Value = '25.877.89'
$ws.Range('E18').Value = '  -0.60%  '

$ws.Range('D19').Value = '''1.001'
$ws.Range('E19').Value = '  -0.03%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''194.82'
$ws.Range('E20').Value = '  -1.30%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''4.413'
$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').Value = '''9.854'
$ws.Range('E22').Value = '  -0.58%  '

$ws.Range('D23').Value = '''6.018'
$ws.Range('E23').Value = '  -0.13%  '

$ws.Range('D24').Value = '''1.003'
$ws.Range('E24').Value = '  -0.02%  '

$ws.Range('D25').Value = '''1.891'
$ws.Range('E25').Value = '  +1.93%  '

$ws.Range('D26').Value = '''141.97'
$ws.Range('E26').Value = '  +0.80%  '

$ws.Range('D27').Value = '''0.1254'
$ws.Range('E27').Value = '  +5.10%  '

$ws.Range('D28').Value = '''6.756'
$ws.Range('E28').Value = '  -0.83%  '

$ws.Range('D29').Value = '''15.55'
$ws.Range('E29').Value = '  -0.20%  '

$ws.Range('D30').Value = '''1.239'
$ws.Range('E30').Value = '  +0.52%  '

$ws.Range('D31').Value = '''0.04897'
$ws.Range('E31').Value = '  +0.96%  '

$ws.Range('D32').Value = '''3.232'
$ws.Range('E32').Value = '  -0.48%  '

$ws.Range('D33').Value = '''3.178'
$ws.Range('E33').Value = '  +0.61%  '

$ws.Range('D34').Value = '''1.545'
$ws.Range('E34').Value = '  +1.51%  '

$ws.Range('D35').Value = '''2.373'
$ws.Range('E35').Value = '  +0.53%  '

$ws.Range('D36').Value = '''0.8943'
$ws.Range('E36').Value = '  -0.29%  '

$ws.Range('D37').Value = '''0.5515'
$ws.Range('E37').Value = '  +1.22%  '

$ws.Range('D38').Value = '''2.538'
$ws.Range('E38').Value = '  -1.56%  '

$ws.Range('D39').Value = '1.114.60'
$ws.Range('E39').Value = '  -2.17%  '

$ws.Range('D40').Value = '''0.01553'
$ws.Range('E40').Value = '  -0.42%  '

$ws.Range('E41').Value = '  -0.12%  '

$ws.Range('D42').Value = '''5.573'
$ws.Range('E42').Value = '  +3.50%  '

$ws.Range('D43').Value = '''0.7937'
$ws.Range('E43').Value = '  -1.93%  '

$ws.Range('D44').Value = '''97.71'
$ws.Range('E44').Value = '  -1.67%  '

$ws.Range('D45').Value = '1.781.06'
$ws.Range('E45').Value = '  +0.08%  '

$ws.Range('D46').Value = '0.0₈114'
$ws.Range('E46').Value = '  -11.29%  '

$ws.Range('D47').Value = '''0.4434'
$ws.Range('E47').Value = '  -1.99%  '

$ws.Range('D48').Value = '''1.002'
$ws.Range('E48').Value = '  +0.34%  '

$ws.Range('D49').Value = '''54.70'
$ws.Range('E49').Value = '  -0.09%  '

$ws.Range('D50').Value = '''0.05130'
$ws.Range('E50').Value = '  +1.53%  '

$ws.Range('D51').Value = '''7.524'
$ws.Range('E51').Value = '  +2.71%  '
